# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.672.44"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.813.99"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "668.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.812.66"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.465"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.04"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.459.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.812.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.587.45"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +16.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.25"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.35"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.965.30"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.56%  "
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.43"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.74"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.178"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.22"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.773.97"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("B38").Value = "Binance-PegBSC-USD"
$ws.Range("C38").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.48"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.01"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.969"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.31%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.56"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.66%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.67"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.301"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.42"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.11%  "
